$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 103.4766596666667
$ws.Range("N2").Value = 310.429979
$ws.Range("O2").Value = 0.877785331764719
$ws.Range("P2").Value = 0.8777853317647188
$ws.Range("Q2").Value = 217.1859192544507
$ws.Range("R2").Value = 1954.673273290056
$ws.Range("S2").Value = 0.0950161944061057
$ws.Range("T2").Value = 0.09501619440610568

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8265796666666668
$ws.Range("N3").Value = 2.479739
$ws.Range("O3").Value = 0.007011818020336602
$ws.Range("P3").Value = 0.0070118180203366
$ws.Range("Q3").Value = 1.734898143410667
$ws.Range("R3").Value = 15.614083290696
$ws.Range("S3").Value = 0.0007589968071363436
$ws.Range("T3").Value = 0.0007589968071363433

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.58054833333333
$ws.Range("N4").Value = 40.741645
$ws.Range("O4").Value = 0.1152028502149446
$ws.Range("P4").Value = 0.1152028502149446
$ws.Range("Q4").Value = 28.50404993025333
$ws.Range("R4").Value = 256.53644937228
$ws.Range("S4").Value = 0.01247017467260965
$ws.Range("T4").Value = 0.01247017467260964

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 103.4766596666667
$ws.Range("N5").Value = 310.429979
$ws.Range("O5").Value = 0.877785331764719
$ws.Range("P5").Value = 0.8777853317647188
$ws.Range("Q5").Value = 1642.706838355105
$ws.Range("R5").Value = 14784.36154519595
$ws.Range("S5").Value = 0.7186642340405289
$ws.Range("T5").Value = 0.7186642340405287

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8265796666666668
$ws.Range("N6").Value = 2.479739
$ws.Range("O6").Value = 0.007011818020336602
$ws.Range("P6").Value = 0.0070118180203366
$ws.Range("Q6").Value = 13.12207096027878
$ws.Range("R6").Value = 118.098638642509
$ws.Range("S6").Value = 0.005740746221728242
$ws.Range("T6").Value = 0.00574074622172824

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.58054833333333
$ws.Range("N7").Value = 40.741645
$ws.Range("O7").Value = 0.1152028502149446
$ws.Range("P7").Value = 0.1152028502149446
$ws.Range("Q7").Value = 215.5931558637772
$ws.Range("R7").Value = 1940.338402773995
$ws.Range("S7").Value = 0.09431937982212775
$ws.Range("T7").Value = 0.09431937982212772

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.4766596666667
$ws.Range("N8").Value = 310.429979
$ws.Range("O8").Value = 0.877785331764719
$ws.Range("P8").Value = 0.8777853317647188
$ws.Range("Q8").Value = 146.5295725942186
$ws.Range("R8").Value = 1318.766153347968
$ws.Range("S8").Value = 0.06410490331808445
$ws.Range("T8").Value = 0.06410490331808444

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8265796666666668
$ws.Range("N9").Value = 2.479739
$ws.Range("O9").Value = 0.007011818020336602
$ws.Range("P9").Value = 0.0070118180203366
$ws.Range("Q9").Value = 1.170489709098667
$ws.Range("R9").Value = 10.534407381888
$ws.Range("S9").Value = 0.0005120749914720172
$ws.Range("T9").Value = 0.000512074991472017

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.58054833333333
$ws.Range("N10").Value = 40.741645
$ws.Range("O10").Value = 0.1152028502149446
$ws.Range("P10").Value = 0.1152028502149446
$ws.Range("Q10").Value = 19.23092559509333
$ws.Range("R10").Value = 173.07833035584
$ws.Range("S10").Value = 0.008413295720207225
$ws.Range("T10").Value = 0.008413295720207221
